$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated activity data per row (columns B-J), reflecting the refreshed
# chart timeframe pulled from the schedule.
$rowData = @{
    2 = @{ "B"=13; "C"=764.322; "D"=1112.327; "E"=66; "F"=12; "G"=34; "H"=262; "I"=2419.3; "J"=-54.02277518290415 }
    3 = @{ "B"=2; "C"=177; "D"=192; "E"=12; "F"=0; "G"=5; "H"=0; "I"=289; "J"=-33.56401384083046 }
    4 = @{ "B"=0; "C"=540; "D"=565; "E"=27; "F"=3; "G"=22; "H"=0; "I"=112; "J"=404.4642857142857 }
    5 = @{ "B"=10; "C"=1041; "D"=1186; "E"=58; "F"=2; "G"=11; "H"=75; "I"=1469; "J"=-19.2648059904697 }
    6 = @{ "B"=6; "C"=930; "D"=976; "E"=38; "F"=6; "G"=11; "H"=0; "I"=1148; "J"=-14.98257839721254 }
    7 = @{ "B"=0; "C"=260; "D"=266; "E"=6; "F"=0; "G"=0; "H"=0; "I"=390; "J"=-31.79487179487179 }
    8 = @{ "B"=0; "C"=98; "D"=131; "E"=33; "F"=0; "G"=2; "H"=0; "I"=124; "J"=5.645161290322576 }
    9 = @{ "B"=0; "C"=425; "D"=430; "E"=2; "F"=2; "G"=12; "H"=1; "I"=1339; "J"=-67.88648244958925 }
    10 = @{ "B"=0; "C"=93; "D"=94; "E"=1; "F"=1; "G"=2; "H"=0; "I"=202; "J"=-53.46534653465347 }
    11 = @{ "B"=0; "C"=4; "D"=4; "E"=0; "F"=0; "G"=0; "H"=0; "I"=6; "J"=-33.33333333333334 }
    12 = @{ "B"=0; "C"=50; "D"=52; "E"=2; "F"=0; "G"=1; "H"=0; "I"=157; "J"=-66.87898089171975 }
}

foreach ($r in $rowData.Keys) {
    $cols = $rowData[$r]
    foreach ($col in $cols.Keys) {
        $ws.Range("$col$r").Value = $cols[$col]
    }
}
